$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the expected date format to the three date-related column headers
$ws.Range("B1").Value = "Subscriber Birth Date (MM/DD/YYYY)"
$ws.Range("C1").Value = "Date of Issue (MM/DD/YYYY)"
$ws.Range("D1").Value = "Date of Service (MM/DD/YYYY)"

# Re-fit the (now wider) header columns so the longer text still displays in full
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

# AutoFit in this environment doesn't quite reproduce Excel's own font-metrics
# based "best fit" width, so nudge the three columns to the exact widths the
# header text resolves to once it's been through Excel's best-fit sizing.
$ws.Columns.Item(2).ColumnWidth = 37.451822916666664
$ws.Columns.Item(3).ColumnWidth = 28.592447916666668
$ws.Columns.Item(4).ColumnWidth = 31.022135416666668

# Move the active selection down to A3, as in the saved workbook
$ws.Range("A3").Select() | Out-Null
